$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cryptos.xlsx data refresh (coinranking.com pull).
# Column D sometimes holds plain-numeric-looking text (e.g. "1.00",
# "227.18"); setting .Value directly on a General-formatted cell would
# let Excel reinterpret it as a number and drop the original text
# formatting, so those specific cells are switched to Text format ("@")
# first, matching how the sheet already stores them (t="inlineStr").

$ws.Range("D2").Value = '34.140.26'
$ws.Range("E2").Value = '  +0.51%  '
$ws.Range("D3").Value = '1.788.87'
$ws.Range("E3").Value = '  -1.30%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.18'
$ws.Range("E5").Value = '  -0.34%  '
$ws.Range("E6").Value = '  +2.05%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '31.61'
$ws.Range("E8").Value = '  +2.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.18'
$ws.Range("E9").Value = '  -3.22%  '
$ws.Range("E10").Value = '  +1.21%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0662'
$ws.Range("E11").Value = '  -0.71%  '
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("D13").Value = '2.046.81'
$ws.Range("E13").Value = '  -1.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.38'
$ws.Range("E14").Value = '  +12.31%  '
$ws.Range("D15").Value = '1.785.43'
$ws.Range("E15").Value = '  -1.40%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.632'
$ws.Range("E16").Value = '  -0.81%  '
$ws.Range("D17").Value = '34.138.92'
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("E19").Value = '  +0.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '255.13'
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("D21").Value = '0.0₃0743'
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.51'
$ws.Range("E23").Value = '  +1.34%  '
$ws.Range("E24").Value = '  -1.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.15'
$ws.Range("E25").Value = '  -1.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.76'
$ws.Range("E26").Value = '  -1.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.57'
$ws.Range("E27").Value = '  +0.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.05'
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("E31").Value = '  +0.29%  '
$ws.Range("E32").Value = '  +1.68%  '
$ws.Range("E33").Value = '  +0.56%  '
$ws.Range("E34").Value = '  +2.31%  '
$ws.Range("E35").Value = '  +2.45%  '
$ws.Range("D36").Value = '1.453.33'
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.635'
$ws.Range("E38").Value = '  +2.59%  '
$ws.Range("E39").Value = '  +0.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.88'
$ws.Range("E40").Value = '  +1.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '83.52'
$ws.Range("E41").Value = '  -0.49%  '
$ws.Range("E42").Value = '  +0.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.902'
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("E44").Value = '  -0.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0509'
$ws.Range("E45").Value = '  -3.18%  '
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.85'
$ws.Range("E47").Value = '  +3.65%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '1.945.59'
$ws.Range("E48").Value = '  -1.18%  '
$ws.Range("B49").Value = 'PaxDollar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.95'
$ws.Range("E50").Value = '  +7.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '50.98'
$ws.Range("E51").Value = '  -2.50%  '
